$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

$ws.Range("E10").Value = "Valo"
$ws.Range("E17").Value = "ValoAikaLeima33"
$ws.Range("E18").Value = "ValoAikaLeima66"
$ws.Range("E19").Value = "ValoAikaLeima100"
$ws.Range("E16").Value = "ValoAikaLeimaONOFF"
$ws.Range("E12").Value = "ValoONOFF"
$ws.Range("E13").Value = "Valo33"
$ws.Range("E13").NumberFormat = "0%"
$ws.Range("E14").Value = "Valo66"
$ws.Range("E15").Value = "Valo100"
$ws.Range("F19").Value = "fk"
$ws.Range("F20").Value = ""
